$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date-serial (D), volume (J), min price (K), max price (L),
# weighted-avg price (M), price $/Kg (P), quality (I)
$rows = @(
    @(512,44637,2000,950,1000,975,975,"Primera"),
    @(513,44637,1240,850,900,875,875,"Segunda"),
    @(514,44208,2300,600,700,650,650,"Primera"),
    @(515,44208,1460,500,550,525,525,"Segunda"),
    @(516,44445,3200,600,700,650,650,"Primera"),
    @(517,44445,1500,500,550,525,525,"Segunda"),
    @(518,44524,3000,600,700,650,650,"Primera"),
    @(519,44524,1600,500,550,525,525,"Segunda"),
    @(520,44355,2400,600,700,650,650,"Primera"),
    @(521,44355,1460,500,550,525,525,"Segunda"),
    @(522,44530,2300,600,650,625,625,"Primera"),
    @(523,44530,1300,500,550,525,525,"Segunda"),
    @(524,44483,2000,600,700,650,650,"Primera"),
    @(525,44483,1320,500,550,525,525,"Segunda"),
    @(526,44294,2400,650,700,675,675,"Primera"),
    @(527,44294,1360,550,600,575,575,"Segunda"),
    @(528,44617,2480,850,900,875,875,"Primera"),
    @(529,44617,1280,750,800,775,775,"Segunda"),
    @(530,44557,2500,600,700,650,650,"Primera"),
    @(531,44557,1480,500,550,525,525,"Segunda"),
    @(532,44489,3000,600,700,650,650,"Primera"),
    @(533,44489,1600,500,550,525,525,"Segunda"),
    @(534,44264,3000,800,900,850,850,"Primera"),
    @(535,44264,1560,700,750,725,725,"Segunda"),
    @(536,44396,3080,700,800,750,750,"Primera"),
    @(537,44396,1400,500,600,550,550,"Segunda"),
    @(538,44232,3000,650,700,675,675,"Primera"),
    @(539,44232,1600,550,600,575,575,"Segunda"),
    @(540,44279,3400,700,750,725,725,"Primera"),
    @(541,44279,1800,600,650,625,625,"Segunda"),
    @(542,44330,3400,650,700,675,675,"Primera"),
    @(543,44330,1680,500,550,525,525,"Segunda"),
    @(544,44504,2000,600,700,650,650,"Primera"),
    @(545,44504,1300,500,550,525,525,"Segunda"),
    @(546,44572,2360,600,700,650,650,"Primera"),
    @(547,44572,1320,500,550,525,525,"Segunda"),
    @(548,44257,2600,850,900,875,875,"Primera"),
    @(549,44257,1500,750,800,775,775,"Segunda"),
    @(550,44301,2500,650,700,675,675,"Primera"),
    @(551,44301,1360,550,600,575,575,"Segunda"),
    @(552,44370,3400,600,700,650,650,"Primera"),
    @(553,44370,1800,500,550,525,525,"Segunda"),
    @(554,44487,2400,600,700,650,650,"Primera"),
    @(555,44487,1400,500,550,525,525,"Segunda"),
    @(556,44174,2800,550,600,575,575,"Primera"),
    @(557,44174,1560,450,500,475,475,"Segunda"),
    @(558,44200,3000,650,700,675,675,"Primera"),
    @(559,44200,1600,550,600,575,575,"Segunda"),
    @(560,44385,2200,600,700,650,650,"Primera"),
    @(561,44385,1300,500,550,525,525,"Segunda"),
    @(562,44236,2400,650,700,675,675,"Primera"),
    @(563,44236,1500,550,600,575,575,"Segunda"),
    @(564,44221,3000,650,700,675,675,"Primera"),
    @(565,44221,1600,550,600,575,575,"Segunda"),
    @(566,44413,2000,650,700,675,675,"Primera"),
    @(567,44413,1400,550,600,575,575,"Segunda"),
    @(568,44272,3400,850,900,875,875,"Primera"),
    @(569,44272,1800,750,800,775,775,"Segunda"),
    @(570,44229,2400,650,700,675,675,"Primera"),
    @(571,44229,1300,550,600,575,575,"Segunda"),
    @(572,44214,3000,600,650,625,625,"Primera"),
    @(573,44214,1600,500,550,525,525,"Segunda"),
    @(574,44299,2600,650,700,675,675,"Primera"),
    @(575,44299,1480,550,600,575,575,"Segunda"),
    @(576,44610,2500,850,900,875,875,"Primera"),
    @(577,44610,1300,750,800,775,775,"Segunda"),
    @(578,44312,3000,650,700,675,675,"Primera"),
    @(579,44312,1480,550,600,575,575,"Segunda"),
    @(580,44399,2200,600,700,650,650,"Primera"),
    @(581,44399,1400,500,550,525,525,"Segunda"),
    @(582,44615,2600,850,900,875,875,"Primera"),
    @(583,44615,1320,750,800,775,775,"Segunda"),
    @(584,44522,2400,600,700,650,650,"Primera"),
    @(585,44522,1500,500,550,525,525,"Segunda"),
    @(586,44543,2500,600,700,650,650,"Primera"),
    @(587,44543,1460,500,550,525,525,"Segunda"),
    @(588,44167,2900,450,500,475,475,"Primera"),
    @(589,44167,1600,350,400,375,375,"Segunda"),
    @(590,44277,2800,700,800,750,750,"Primera"),
    @(591,44277,1400,600,650,625,625,"Segunda"),
    @(592,44258,3400,850,900,875,875,"Primera"),
    @(593,44258,2000,750,800,775,775,"Segunda"),
    @(594,44390,2000,600,700,650,650,"Primera"),
    @(595,44390,1400,500,550,525,525,"Segunda"),
    @(596,44349,3360,600,700,650,650,"Primera"),
    @(597,44349,1800,500,550,525,525,"Segunda"),
    @(598,44285,2500,650,700,675,675,"Primera"),
    @(599,44285,1460,550,600,575,575,"Segunda"),
    @(600,44498,3360,600,700,650,650,"Primera"),
    @(601,44498,1600,500,550,525,525,"Segunda"),
    @(602,44179,2500,550,600,575,575,"Primera"),
    @(603,44179,1560,450,500,475,475,"Segunda"),
    @(604,44418,2000,650,700,675,675,"Primera"),
    @(605,44418,1400,550,600,575,575,"Segunda"),
    @(606,44595,2000,750,800,775,775,"Primera"),
    @(607,44595,1160,650,700,675,675,"Segunda"),
    @(608,44628,2400,950,1000,975,975,"Primera"),
    @(609,44628,1300,850,900,875,875,"Segunda"),
    @(610,44335,3400,650,700,675,675,"Primera"),
    @(611,44335,1800,500,550,525,525,"Segunda"),
    @(612,44552,3200,600,700,650,650,"Primera"),
    @(613,44552,1600,500,550,525,525,"Segunda"),
    @(614,44544,2400,600,650,625,625,"Primera"),
    @(615,44544,1320,500,550,525,525,"Segunda"),
    @(616,44160,2800,450,500,475,475,"Primera"),
    @(617,44160,1600,350,400,375,375,"Segunda")
)

foreach ($row in $rows) {
    $r = $row[0]
    $dateVal = $row[1]
    $volVal = $row[2]
    $minVal = $row[3]
    $maxVal = $row[4]
    $avgVal = $row[5]
    $pkgVal = $row[6]
    $quality = $row[7]

    if ($r -ge 616) {
        # Brand-new rows appended at the bottom of the table: populate every
        # column, matching the constant values used throughout the sheet.
        $ws.Cells.Item($r, 1).Value = 8
        $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
        $ws.Cells.Item($r, 3).Value = "Coquimbo"
        $ws.Cells.Item($r, 5).Value = 4
        $ws.Cells.Item($r, 6).Value = 100112023
        $ws.Cells.Item($r, 7).Value = "Brócoli"
        $ws.Cells.Item($r, 8).Value = "Sin especificar"
        $ws.Cells.Item($r, 9).Value = $quality
        $ws.Cells.Item($r, 14).Value = "$/unidad"
        $ws.Cells.Item($r, 15).Value = "Provincia del Elquí"
        $ws.Cells.Item($r, 17).Value = 1
        $ws.Cells.Item($r, 18).Value = "Hortaliza"
    }

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $dateVal
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 10).Value = $volVal
    $ws.Cells.Item($r, 11).Value = $minVal
    $ws.Cells.Item($r, 12).Value = $maxVal
    $ws.Cells.Item($r, 13).Value = $avgVal
    $ws.Cells.Item($r, 16).Value = $pkgVal
}
